# MHD2-259: Report template and related changes for reporting on 136 genes
#
# The document body change is a single table-cell shading colour update:
# the title/header cell's background fill moves from E8E7EC to ECEAF2.

$d = $word.ActiveDocument

# The document contains exactly one table, whose single cell holds the
# whole clinical-context body; its header shading is what changes colour.
$cell = $d.Tables(1).Cell(1, 1)

# Word COM colour longs are packed as 0x00BBGGRR (BGR), i.e. Red + Green*256 + Blue*65536.
$newFill = 0xEC + (0xEA * 0x100) + (0xF2 * 0x10000)
$cell.Shading.BackgroundPatternColor = $newFill
